# Correction in sa algorithm and 746 logs
# Update the Fitness (column C) values for rows 2-252 according to the
# corrected run_28 log data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2..44 (Generation 0..42) -> 7312
$ws.Range("C2:C44").Value = 7312

# Rows 45..168 (Generation 43..166) -> 7295
$ws.Range("C45:C168").Value = 7295

# Rows 169..252 (Generation 167..250) -> 7293
$ws.Range("C169:C252").Value = 7293
